# Rename the "_old"/"_new" header-column suffixes to the respective
# input-file format-version suffixes ("_FV2304" / "_FV2310"), freeze the
# header row, and format the header range as an Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2304", "Segmentgruppe_FV2304", "Segment_FV2304", "Datenelement_FV2304", "Segment ID_FV2304",
  "Code_FV2304", "Qualifier_FV2304", "Beschreibung_FV2304", "Bedingungsausdruck_FV2304", "Bedingung_FV2304",
  "diff",
  "Segmentname_FV2310", "Segmentgruppe_FV2310", "Segment_FV2310", "Datenelement_FV2310", "Segment ID_FV2310",
  "Code_FV2310", "Qualifier_FV2310", "Beschreibung_FV2310", "Bedingungsausdruck_FV2310", "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $col = $i + 1
  $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Freeze the header row (top row) so it stays visible while scrolling.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the header + data range into a real Excel Table ("Table1").
$tableRange = $ws.Range("A1:U64")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""
